$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmark")
$ws.Columns.Item(3).Insert()
$ws.Cells.Item(1, 3).Value = "software_backend"
